$wb = $excel.ActiveWorkbook

# 1. Create the new "Slovakia" sheet by copying the "Portugal" sheet (same
#    column widths / styles / merged cells / page setup as the other
#    per-market tabs) and placing it immediately after Portugal.
$portugal = $wb.Worksheets.Item("Portugal")
$portugal.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

# 2. Market name / NGC reference for the new tab.
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3219/T3183/T3238 "

# 3. Replace the product list (rows 8-32) with the Slovakia line-up.
$products = @("Pro215S","Pro215D","Pro815D","Pro885D","Pro16xD","Pro16xBB","Pro32xD","Pro32xBB","P115S","P115D","P405D","P485D","P805D","P885D","MZX125","MZX250","MZX251","MZX252","MZX253","MZX254","ZX1","ZX4","Black Box","Wg","Panels")
$row = 8
foreach ($product in $products) {
    $slovakia.Range("A" + $row).Value = $product
    $row = $row + 1
}

# The copied sheet still has the four extra (now unused) rows the Portugal
# tab had at the bottom of its list - drop them so the sheet ends at row 32.
$slovakia.Rows.Item(33).Delete()
$slovakia.Rows.Item(33).Delete()
$slovakia.Rows.Item(33).Delete()
$slovakia.Rows.Item(33).Delete()

# 4. Update the selections left on the other tabs as part of this review.
$belgium = $wb.Worksheets.Item("Belgium")
$belgium.Range("A33:A34").Select()

$portugal.Range("A33").Select()

# 5. Finally leave the new Slovakia sheet selected/active.
$slovakia.Range("B18").Select()
